$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 2.989631
$ws.Range("H2").Value = 5.979262
$ws.Range("I2").Value = 0.09195719396405833
$ws.Range("J2").Value = 0.06346072881692182
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.704179666666667
$ws.Range("N2").Value = 5.112539
$ws.Range("O2").Value = 0.3144996488703566
$ws.Range("P2").Value = 0.3144996488703566
$ws.Range("Q2").Value = 5.094868361036333
$ws.Range("R2").Value = 30.569210166218
$ws.Range("S2").Value = 0.02892050521279962
$ws.Range("T2").Value = 0.01995837692997883

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 2.989631
$ws.Range("H3").Value = 5.979262
$ws.Range("I3").Value = 0.09195719396405833
$ws.Range("J3").Value = 0.06346072881692182
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.714521666666666
$ws.Range("N3").Value = 11.143565
$ws.Range("O3").Value = 0.6855003511296434
$ws.Range("P3").Value = 0.6855003511296432
$ws.Range("Q3").Value = 11.10504912483833
$ws.Range("R3").Value = 66.63029474903
$ws.Range("S3").Value = 0.0630366887512587
$ws.Range("T3").Value = 0.04350235188694297

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.926294333333334
$ws.Range("H4").Value = 5.778883
$ws.Range("I4").Value = 0.05925032943604069
$ws.Range("J4").Value = 0.06133401194457102
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.704179666666667
$ws.Range("N4").Value = 5.112539
$ws.Range("O4").Value = 0.3144996488703566
$ws.Range("P4").Value = 0.3144996488703566
$ws.Range("Q4").Value = 3.282751634881889
$ws.Range("R4").Value = 29.544764713937
$ws.Range("S4").Value = 0.01863420780308776
$ws.Range("T4").Value = 0.01928952522037785

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.926294333333334
$ws.Range("H5").Value = 5.778883
$ws.Range("I5").Value = 0.05925032943604069
$ws.Range("J5").Value = 0.06133401194457102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.714521666666666
$ws.Range("N5").Value = 11.143565
$ws.Range("O5").Value = 0.6855003511296434
$ws.Range("P5").Value = 0.6855003511296432
$ws.Range("Q5").Value = 7.155262037543889
$ws.Range("R5").Value = 64.397358337895
$ws.Range("S5").Value = 0.04061612163295294
$ws.Range("T5").Value = 0.04204448672419317

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.330727000000001
$ws.Range("H6").Value = 24.992181
$ws.Range("I6").Value = 0.256242418746868
$ws.Range("J6").Value = 0.2652538090795195
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.704179666666667
$ws.Range("N6").Value = 5.112539
$ws.Range("O6").Value = 0.3144996488703566
$ws.Range("P6").Value = 0.3144996488703566
$ws.Range("Q6").Value = 14.197055561951
$ws.Range("R6").Value = 127.773500057559
$ws.Range("S6").Value = 0.08058815072158089
$ws.Range("T6").Value = 0.0834222298170335

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.330727000000001
$ws.Range("H7").Value = 24.992181
$ws.Range("I7").Value = 0.256242418746868
$ws.Range("J7").Value = 0.2652538090795195
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.714521666666666
$ws.Range("N7").Value = 11.143565
$ws.Range("O7").Value = 0.6855003511296434
$ws.Range("P7").Value = 0.6855003511296432
$ws.Range("Q7").Value = 30.944665940585
$ws.Range("R7").Value = 278.501993465265
$ws.Range("S7").Value = 0.1756542680252871
$ws.Range("T7").Value = 0.181831579262486

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.84396866666667
$ws.Range("H8").Value = 47.531906
$ws.Range("I8").Value = 0.4873400429153729
$ws.Range("J8").Value = 0.5044785454822717
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.704179666666667
$ws.Range("N8").Value = 5.112539
$ws.Range("O8").Value = 0.3144996488703566
$ws.Range("P8").Value = 0.3144996488703566
$ws.Range("Q8").Value = 27.00096924103711
$ws.Range("R8").Value = 243.008723169334
$ws.Range("S8").Value = 0.1532682723773493
$ws.Range("T8").Value = 0.1586583254168027

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.84396866666667
$ws.Range("H9").Value = 47.531906
$ws.Range("I9").Value = 0.4873400429153729
$ws.Range("J9").Value = 0.5044785454822717
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.714521666666666
$ws.Range("N9").Value = 11.143565
$ws.Range("O9").Value = 0.6855003511296434
$ws.Range("P9").Value = 0.6855003511296432
$ws.Range("Q9").Value = 58.85276489832111
$ws.Range("R9").Value = 529.67488408489
$ws.Range("S9").Value = 0.3340717705380236
$ws.Range("T9").Value = 0.3458202200654689

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.096654
$ws.Range("H10").Value = 9.289962000000001
$ws.Range("I10").Value = 0.09524908342119047
$ws.Range("J10").Value = 0.09859875001321379
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.704179666666667
$ws.Range("N10").Value = 5.112539
$ws.Range("O10").Value = 0.3144996488703566
$ws.Range("P10").Value = 0.3144996488703566
$ws.Range("Q10").Value = 5.277254781502001
$ws.Range("R10").Value = 47.495293033518
$ws.Range("S10").Value = 0.02995580329118771
$ws.Range("T10").Value = 0.03100927225821181

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.096654
$ws.Range("H11").Value = 9.289962000000001
$ws.Range("I11").Value = 0.09524908342119047
$ws.Range("J11").Value = 0.09859875001321379
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.714521666666666
$ws.Range("N11").Value = 11.143565
$ws.Range("O11").Value = 0.6855003511296434
$ws.Range("P11").Value = 0.6855003511296432
$ws.Range("Q11").Value = 11.50258837717
$ws.Range("R11").Value = 103.52329539453
$ws.Range("S11").Value = 0.06529328013000275
$ws.Range("T11").Value = 0.06758947775500197

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.323841
$ws.Range("H12").Value = 0.647682
$ws.Range("I12").Value = 0.009960931516469628
$ws.Range("J12").Value = 0.00687415466350221
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.704179666666667
$ws.Range("N12").Value = 5.112539
$ws.Range("O12").Value = 0.3144996488703566
$ws.Range("P12").Value = 0.3144996488703566
$ws.Range("Q12").Value = 0.551883247433
$ws.Range("R12").Value = 3.311299484598
$ws.Range("S12").Value = 0.003132709464351367
$ws.Range("T12").Value = 0.00216191922795197

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.323841
$ws.Range("H13").Value = 0.647682
$ws.Range("I13").Value = 0.009960931516469628
$ws.Range("J13").Value = 0.00687415466350221
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.714521666666666
$ws.Range("N13").Value = 11.143565
$ws.Range("O13").Value = 0.6855003511296434
$ws.Range("P13").Value = 0.6855003511296432
$ws.Range("Q13").Value = 1.202914411055
$ws.Range("R13").Value = 7.217486466329999
$ws.Range("S13").Value = 0.006828222052118261
$ws.Range("T13").Value = 0.004712235435550239

Write-Host "Updated H2-M3/Klrd1 LR-pair values (natmi rerun per Dr Hou advice)"